$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells
$ws.Range("G1").Value = "Elapsed Time"
$ws.Range("H1").Value = "CPU"

# Copy header formatting (bold font, thin border, centered alignment) from A1
$ws.Range("A1").Copy()
$ws.Range("G1:H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill Elapsed Time / CPU columns for rows 2-10
for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 7).Value = 1.127317944850074
    $ws.Cells.Item($r, 8).Value = 0.985
}

# Update slightly changed MSE/R2/MAE values (B,C,D) per row
$ws.Range("B2").Value = 0.4252561148737903
$ws.Range("C2").Value = 0.930721389331257
$ws.Range("D2").Value = 0.5068547933339305

$ws.Range("B3").Value = 2.878893579773744
$ws.Range("C3").Value = 0.9588727023492972
$ws.Range("D3").Value = 1.266131410638987

$ws.Range("B4").Value = 1.09214980763019
$ws.Range("C4").Value = 0.9460522249730414
$ws.Range("D4").Value = 0.8035090531811903

$ws.Range("B5").Value = 1.56219342376502
$ws.Range("D5").Value = 0.9389581486305347

$ws.Range("B7").Value = 1.271977756333027
$ws.Range("D7").Value = 0.8344729122263086

$ws.Range("B9").Value = 4.893487171018756
$ws.Range("C9").Value = 0.9415017423165769
$ws.Range("D9").Value = 1.680542520641608

$ws.Range("B10").Value = 0.8506165638197862
$ws.Range("C10").Value = 0.9974643383697318
$ws.Range("D10").Value = 0.7293641098176615
